# "formateado a 2 decimales"
# Insert a bold header row ("Puntaje" / "Datos del vino") and reformat the
# ranking scores in column A from raw float strings (e.g. 92.66666666666667)
# to comma-decimal, 2-decimal-place strings (e.g. 92,67).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 10 data rows down to make room for a header row.
$ws.Rows("1:1").Insert()

# New header row.
$ws.Range("A1").Value = "Puntaje"
$ws.Range("B1").Value = "Datos del vino"
$ws.Range("A1:B1").Font.Bold = $true

# Re-write the wine names (unchanged text) together with the newly
# formatted scores so every row ends up consistent.
$ws.Range("A2").Value = "92,67"
$ws.Range("B2").Value = "Montes Alpha Cabernet Sauvignon - 4500.0 - Bodega Montes, Colchagua Andes, Colchagua, Chile - [Cabernet Gran Reserva]"

$ws.Range("A3").Value = "92,33"
$ws.Range("B3").Value = "Norton Chardonnay Reserva - 3100.0 - Bodega Norton, Valle de Uco, Mendoza, Argentina - [Chardonnay Clásico]"

$ws.Range("A4").Value = "91,33"
$ws.Range("B4").Value = "Garzón Albariño - 3800.0 - Bodega Garzón, Las Violetas, Canelones, Uruguay - [Sauvignon Blanc Finca]"

$ws.Range("A5").Value = "91,33"
$ws.Range("B5").Value = "Norton Sauvignon Blanc - 3000.0 - Bodega Norton, Valle de Uco, Mendoza, Argentina - [Sauvignon Blanc Finca]"

$ws.Range("A6").Value = "91,00"
$ws.Range("B6").Value = "Norton Malbec - 2900.0 - Bodega Norton, Valle de Uco, Mendoza, Argentina - [Malbec Reserva]"

$ws.Range("A7").Value = "91,00"
$ws.Range("B7").Value = "Concha y Toro Chardonnay - 3400.0 - Bodega Concha y Toro, Maipo Alto, Maipo, Chile - [Chardonnay Clásico]"

$ws.Range("A8").Value = "91,00"
$ws.Range("B8").Value = "Garzón Tannat - 4000.0 - Bodega Garzón, Las Violetas, Canelones, Uruguay - [Pinot Noir Reserva]"

$ws.Range("A9").Value = "90,00"
$ws.Range("B9").Value = "Trapiche Malbec - 2500.0 - Bodega Trapiche, Luján de Cuyo, Mendoza, Argentina - [Malbec Reserva]"

$ws.Range("A10").Value = "90,00"
$ws.Range("B10").Value = "Trapiche Cabernet Sauvignon - 2700.0 - Bodega Trapiche, Luján de Cuyo, Mendoza, Argentina - [Cabernet Gran Reserva]"

$ws.Range("A11").Value = "90,00"
$ws.Range("B11").Value = "Callia Malbec - 2000.0 - Bodega Callia, Tulum, San Juan, Argentina - [Malbec Reserva]"
